$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update F column "想去人数" (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 95
$wsExhibit.Range("F5").Value = 160
$wsExhibit.Range("F6").Value = 137
$wsExhibit.Range("F9").Value = 2062
$wsExhibit.Range("F11").Value = 4945
$wsExhibit.Range("F12").Value = 100

# Sheet "全部类型" (All Types) - update F column "想去人数" (want-to-go count) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 95
$wsAll.Range("F7").Value = 160
$wsAll.Range("F8").Value = 137
$wsAll.Range("F13").Value = 2062
$wsAll.Range("F15").Value = 4946
$wsAll.Range("F16").Value = 100
